$d = $word.ActiveDocument

$find = $d.Content.Find
$find.Execute(
    "github.com/fatihgulsen/FatihGulsen_g171210070_Tasarim",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "github.com/fatihgulsen/FatihGulsen_g171210070_Tasarim",
    2
)
